$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the text cells that are getting brand-new content first, so the old
# shared strings ("Tanah Lapangan 1/2", "Jalan Gerilya 1/2", "Tempat Olahraga
# 1/2") become unreferenced and get dropped from the shared-string table on
# save, instead of lingering as dead/unused entries.
$ws.Cells.Item(2, 3).ClearContents()
$ws.Cells.Item(2, 6).ClearContents()
$ws.Cells.Item(2, 10).ClearContents()
$ws.Cells.Item(3, 3).ClearContents()
$ws.Cells.Item(3, 6).ClearContents()
$ws.Cells.Item(3, 10).ClearContents()

# --- Row 2 (kode_aset 02.01.0001) -------------------------------------------------
$ws.Cells.Item(2, 6).Value = "Jl. Prof. Dr. Suharso No 58 Purwokerto"
$ws.Cells.Item(2, 10).Value = "Kantor Sekretariat"

# --- Row 3 (kode_aset 02.01.0002) -------------------------------------------------
$ws.Cells.Item(3, 3).Value = "Tanah Kwarcab 2"

$ws.Cells.Item(2, 3).Value = "Tanah Kwarcab 1"

$ws.Cells.Item(3, 6).Value = "Desa Kaliori Kec. Kalibagor"
$ws.Cells.Item(3, 10).Value = "Kantor Pusdiklat Pramuka dan Bumi Perkemahan"

# --- Numeric fields ----------------------------------------------------------------
$ws.Cells.Item(2, 5).Value = 5114
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(3, 11).Value = 564100000

# --- Rows 4-12: only the "harga" (K) column resets to 0 ----------------------------
for ($r = 4; $r -le 12; $r++) {
    $ws.Cells.Item($r, 11).Value = 0
}

# --- Column K width: widen from 7 to 12 (still best-fit/custom) --------------------
$ws.Columns("K").ColumnWidth = 11.1428

# --- Selection / scroll position ----------------------------------------------------
$ws.Range("G7").Select()
